$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (ano 2025) with refreshed recurrence metrics
$ws.Range("C8").Value = 966
$ws.Range("D8").Value = 163
$ws.Range("E8").Value = 803
$ws.Range("F8").Value = 6.685808039376538
$ws.Range("G8").Value = 83.12629399585921
$ws.Range("H8").Value = 16.87370600414079
